$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = "https://www.loretech.io/products/teledyne-flir-ffy-u3-04s2m-s?variant=41921567195312&currency=USD&utm_medium=product_sync&utm_source=google&utm_content=sag_organic&utm_campaign=sag_organic&gclid=CjwKCAjwvfmoBhAwEiwAG2tqzJ1Q4EKml1YO4GNgePSl51jbAo5Rcsx44ZdUebNdtIqEfOWRUJiS6hoCbBUQAvD_BwE"
$ws.Range("A7").Value = "Teledyne FLIR FFY-U3-04S2M-S"
$ws.Range("C7").Value = 1

$ws.Range("C14").Select()
